$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (some look like plain decimals, e.g. "1.009").
# Force the column to Text format first so Excel does not silently coerce
# these numeric-looking strings into real numbers when we set .Value below.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.870.63"
$ws.Range("E2").Value = "  -2.19%  "
$ws.Range("D3").Value = "1.778.87"
$ws.Range("E3").Value = "  -2.80%  "
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("D5").Value = "1.007"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "307.86"
$ws.Range("E6").Value = "  -1.52%  "
$ws.Range("D7").Value = "0.4223"
$ws.Range("E7").Value = "  -1.69%  "
$ws.Range("D8").Value = "0.3617"
$ws.Range("E8").Value = "  -1.05%  "
$ws.Range("D9").Value = "0.07190"
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("D10").Value = "0.8381"
$ws.Range("E10").Value = "  -3.23%  "
$ws.Range("D11").Value = "20.26"
$ws.Range("E11").Value = "  -1.89%  "
$ws.Range("D12").Value = "1.801.13"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").Value = "5.252"
$ws.Range("E13").Value = "  -2.74%  "
$ws.Range("D14").Value = "6.336"
$ws.Range("E14").Value = "  -2.87%  "
$ws.Range("D15").Value = "0.06808"
$ws.Range("E15").Value = "  -1.75%  "
$ws.Range("D16").Value = "1.011"
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").Value = "79.17"
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("D18").Value = "0.000008682"
$ws.Range("E18").Value = "  -2.30%  "
$ws.Range("D19").Value = "1.009"
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("D20").Value = "14.95"
$ws.Range("E20").Value = "  -2.92%  "
$ws.Range("D21").Value = "26.972.41"
$ws.Range("E21").Value = "  -2.72%  "
$ws.Range("D22").Value = "5.010"
$ws.Range("E22").Value = "  -2.36%  "
$ws.Range("D23").Value = "11.04"
$ws.Range("E23").Value = "  +2.02%  "
$ws.Range("D24").Value = "2.026.38"
$ws.Range("E24").Value = "  -1.03%  "
$ws.Range("D25").Value = "1.928"
$ws.Range("E25").Value = "  -2.80%  "
$ws.Range("D26").Value = "153.31"
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("D27").Value = "18.14"
$ws.Range("E27").Value = "  -4.18%  "
$ws.Range("D28").Value = "5.041"
$ws.Range("E28").Value = "  -1.88%  "
$ws.Range("D29").Value = "114.23"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "1.629"
$ws.Range("E30").Value = "  -11.30%  "
$ws.Range("D31").Value = "0.08944"
$ws.Range("D32").Value = "0.7202"
$ws.Range("E32").Value = "  -4.55%  "
$ws.Range("D33").Value = "2.845"
$ws.Range("E33").Value = "  -4.40%  "
$ws.Range("D34").Value = "4.328"
$ws.Range("E34").Value = "  -4.69%  "
$ws.Range("D35").Value = "1.091"
$ws.Range("E35").Value = "  -3.84%  "
$ws.Range("D36").Value = "1.008"
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("E37").Value = "  -0.84%  "
$ws.Range("D38").Value = "0.01892"
$ws.Range("E38").Value = "  -2.32%  "
$ws.Range("D39").Value = "0.05079"
$ws.Range("E39").Value = "  -4.58%  "
$ws.Range("D40").Value = "0.1610"
$ws.Range("E40").Value = "  -3.32%  "
$ws.Range("D41").Value = "0.4914"
$ws.Range("E41").Value = "  -3.30%  "
$ws.Range("D42").Value = "2.513"
$ws.Range("E42").Value = "  -10.39%  "
$ws.Range("D43").Value = "6.096"
$ws.Range("E43").Value = "  -7.28%  "
$ws.Range("D44").Value = "7.933"
$ws.Range("E44").Value = "  -4.96%  "
$ws.Range("D45").Value = "1.007"
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("D46").Value = "104.39"
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("D47").Value = "10.05"
$ws.Range("E47").Value = "  -3.60%  "
$ws.Range("D48").Value = "0.06228"
$ws.Range("E48").Value = "  -4.19%  "
$ws.Range("D49").Value = "0.4474"
$ws.Range("E49").Value = "  -4.50%  "
$ws.Range("D50").Value = "1.571"
$ws.Range("E50").Value = "  -2.63%  "
$ws.Range("D51").Value = "1.714"
$ws.Range("E51").Value = "  -1.59%  "
